$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '37.138.82'
$ws.Range('E2').Value2 = '  -1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '2.022.57'
$ws.Range('E3').Value2 = '  -3.16%  '
$ws.Range('E4').Value2 = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '226.97'
$ws.Range('E5').Value2 = '  -2.90%  '
$ws.Range('E6').Value2 = '  -3.15%  '
$ws.Range('E7').Value2 = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '55.09'
$ws.Range('E8').Value2 = '  -5.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.381'
$ws.Range('E9').Value2 = '  -3.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.0790'
$ws.Range('E10').Value2 = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.105'
$ws.Range('E11').Value2 = '  -3.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '2.326.17'
$ws.Range('E12').Value2 = '  -2.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '14.27'
$ws.Range('E13').Value2 = '  -6.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '20.32'
$ws.Range('E14').Value2 = '  -3.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.743'
$ws.Range('E15').Value2 = '  -4.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '5.17'
$ws.Range('E16').Value2 = '  -3.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '2.016.19'
$ws.Range('E17').Value2 = '  -3.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '37.028.00'
$ws.Range('E18').Value2 = '  -1.95%  '
$ws.Range('E19').Value2 = '  -1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '68.87'
$ws.Range('E20').Value2 = '  -3.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '0.0₃0834'
$ws.Range('E21').Value2 = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '223.23'
$ws.Range('E22').Value2 = '  -2.81%  '
$ws.Range('E23').Value2 = '  +0.18%  '
$ws.Range('E24').Value2 = '  -0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '2.27'
$ws.Range('E25').Value2 = '  -5.37%  '
$ws.Range('B26').Value2 = 'Monero'
$ws.Range('C26').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '167.91'
$ws.Range('E26').Value2 = '  -1.96%  '
$ws.Range('B27').Value2 = 'Cosmos'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '9.36'
$ws.Range('E27').Value2 = '  -3.95%  '
$ws.Range('E28').Value2 = '  -7.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '18.71'
$ws.Range('E29').Value2 = '  -4.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '1.33'
$ws.Range('E30').Value2 = '  -4.54%  '
$ws.Range('E31').Value2 = '  -4.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '4.48'
$ws.Range('E32').Value2 = '  -4.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '0.0607'
$ws.Range('E33').Value2 = '  -4.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '4.47'
$ws.Range('E34').Value2 = '  -2.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '2.35'
$ws.Range('E35').Value2 = '  -5.68%  '
$ws.Range('E36').Value2 = '  +0.18%  '
$ws.Range('E37').Value2 = '  +0.18%  '
$ws.Range('E38').Value2 = '  -5.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '5.35'
$ws.Range('E39').Value2 = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '1.492.81'
$ws.Range('E40').Value2 = '  +2.80%  '
$ws.Range('E41').Value2 = '  -7.70%  '
$ws.Range('E42').Value2 = '  -2.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.0928'
$ws.Range('E43').Value2 = '  -4.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '94.82'
$ws.Range('E44').Value2 = '  -6.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '16.53'
$ws.Range('E45').Value2 = '  -0.76%  '
$ws.Range('E46').Value2 = '  -5.64%  '
$ws.Range('E47').Value2 = '  -5.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '7.13'
$ws.Range('E48').Value2 = '  -1.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '2.91'
$ws.Range('E49').Value2 = '  -1.77%  '
$ws.Range('B50').Value2 = 'FTXToken'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '3.69'
$ws.Range('E50').Value2 = '  -10.70%  '
$ws.Range('B51').Value2 = 'RocketPoolETH'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '2.215.47'
$ws.Range('E51').Value2 = '  -2.82%  '
